# Append a newly scraped job listing at the top of the list (row 3) and
# refresh the "fetched at" timestamp on every existing data row.
#
# Sheet layout ("ランサーズ"):
#   row 1       headers
#   row 2..n    one job posting per row (A:取得日時, B:タイトル, C:カテゴリ,
#               D:価格, E:締切, F:URL (hyperlink), G:優先度スコア, H:スキル概要)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2025-10-31 12:37:08"

# Widen the title column slightly to fit the new entry. The engine's
# ColumnWidth -> stored-width conversion adds 5/6 (0.8333...), so back that
# off here to land on an exact stored width of 47.
$ws.Columns.Item(2).ColumnWidth = 46.166666666666664

# Hyperlink objects track their own target independently of the cell text,
# and per-item Delete() is a no-op here, so drop the whole collection before
# shifting rows around and rebuild it afterwards against the final layout.
$ws.Hyperlinks.Delete()

# Make room for the new listing right under the first (unchanged) row.
$ws.Rows.Item(3).Insert()

# Refresh the capture timestamp on every row that was already present.
$ws.Range("A2").Value = $timestamp
$ws.Range("A4").Value = $timestamp
$ws.Range("A5").Value = $timestamp

# New job listing, inserted as the new row 3.
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【急募】Google Workspace/LLM連携!AI診断レポート自動生成システム構築"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5424558"
$ws.Range("G3").Value = 325
$ws.Range("H3").Value = "🔥AI,Ai"

# Re-create the hyperlinks for every URL cell in the (now 4-row) list.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5424032")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5424558")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5422936")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5424258")

# Hyperlinks.Add() re-styles the cell via a freshly-minted (but otherwise
# identical) style record; reapply the named "Hyperlink" style so the URL
# cells keep using the workbook's existing style slot.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
